$wb = $excel.ActiveWorkbook

# Worksheet involved in this edit
$survey = $wb.Worksheets.Item("survey")

# --- "survey" sheet data edits ---
# Row 2: B2 ("admin_name") is cleared out entirely (value removed, style kept);
# C2's prompt text changes from "What is your name?" to "placeholder text".
# Set this new shared string first so it gets the lower of the two new
# shared-string indices (matches the canonical ordering of the saved file).
$survey.Range("C2").Value = "placeholder text"
$survey.Range("B2").ClearContents()

# Row 1: A1 label changes from "type" to "note".
$survey.Range("A1").Value = "note"

# Row 3: A3/B3 ("string"/"notes") are removed completely; C3's text
# ("Record any relevant notes on this distribution") is cleared but keeps
# its style.
$survey.Range("A3").Clear()
$survey.Range("B3").Clear()
$survey.Range("C3").ClearContents()

# --- Active sheet / selection bookkeeping ---
# "survey" becomes the active (selected) tab, with B5 selected; "settings"
# loses its previous tabSelected flag automatically once survey is activated.
$survey.Activate()
$survey.Range("B5").Select()
